$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that are no longer part of the final plan
# (the "Spikes: Square Root Test" duplicate row and the "MakeFile for each Spike" row)
$ws.Range("A5:G6").EntireRow.Delete()

# Rename the remaining tasks to reflect the final scope of work
$ws.Range("A2").Value = "Mock.c Creation"
$ws.Range("A3").Value = "Test script for discriminant test"
$ws.Range("A4").Value = "Absolute/Relative Error Test Scripts"

# Row 3's estimated time moved from 1 hour to 45 min.
$ws.Range("B3").Value = "45 min."

# Code Complete is back down to 0% for all remaining tasks
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0

# Widen the Task column now that the labels are longer (compensate for the
# host's implicit +5/6 character padding so the saved width lands on 31)
$ws.Columns("A").ColumnWidth = 31 - 5/6

# Leave the selection where the author left it when finishing up
$ws.Range("A14").Select() | Out-Null
